# Generate Report for Handback
# - Flips status from "Ready for handoff" to "Handed back: in sync with en-US"
#   (shared text, so every cell that showed the old status shows the new one)
# - Stamps the zh-cn handback datetime (previously the zero date) and gives
#   de-de its own (later) handback datetime
# - Adds "Latest Target File" (F) / "Latest Handback File" (G) columns on the
#   zh-cn and de-de sheets, pointing at the same files that were handed off
#   (since the round trip is complete and nothing changed in translation)

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$zhHandback = "2016-03-24 13:03:10"
$deHandback = "2016-03-24 13:03:24"

$mdUuid1 = "0c5cca4e-0f4f-4441-8e7e-6d7d33aed0f9.md"
$mdUuid2 = "d40e0c57-009a-4ab9-a6f2-5b94d214dfa8.md"
$xlfZh1 = "0c5cca4e-0f4f-4441-8e7e-6d7d33aed0f9.5deef83f6b001d65d7172a080297fd75af593ef5.zh-cn.xlf"
$xlfZh2 = "d40e0c57-009a-4ab9-a6f2-5b94d214dfa8.becc88bdf64d88202c8b38d244f2c24d590775cd.zh-cn.xlf"
$xlfDe1 = "0c5cca4e-0f4f-4441-8e7e-6d7d33aed0f9.5deef83f6b001d65d7172a080297fd75af593ef5.de-de.xlf"
$xlfDe2 = "d40e0c57-009a-4ab9-a6f2-5b94d214dfa8.becc88bdf64d88202c8b38d244f2c24d590775cd.de-de.xlf"

$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/8d0ff1e955987165639646e6a2c8bd8411d786e2/e2e/$mdUuid1"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/8d0ff1e955987165639646e6a2c8bd8411d786e2/e2e/$mdUuid2"
$xlfZh1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/00eb5d7494c1fd4b35bd11be63b9f0df4c407363/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh1"
$xlfZh2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/00eb5d7494c1fd4b35bd11be63b9f0df4c407363/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh2"
$xlfDe1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/72a0940082dcd3bde6f7a882fd51121b495d2124/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe1"
$xlfDe2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/72a0940082dcd3bde6f7a882fd51121b495d2124/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe2"

# ---------------------------------------------------------------------------
# 1. Flip every "Status" cell (Overview!B/C, zh-cn!C, de-de!C) from
#    "Ready for handoff" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column H): zh-cn gets the earlier stamp,
#    de-de gets the later one
# ---------------------------------------------------------------------------
$wsZh.Range("H2").Value = $zhHandback
$wsZh.Range("H3").Value = $zhHandback

$wsDe.Range("H2").Value = $deHandback
$wsDe.Range("H3").Value = $deHandback

# ---------------------------------------------------------------------------
# 3. New columns F (Latest Target File) / G (Latest Handback File) on the
#    zh-cn and de-de sheets, each row linking to the same md/xlf pair used
#    for the original handoff (columns A/D)
# ---------------------------------------------------------------------------

# zh-cn, row 2 (0c5cca4e...)
$wsZh.Range("F2").Value = $mdUuid1
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $md1Url) | Out-Null

$wsZh.Range("G2").Value = $xlfZh1
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $xlfZh1Url) | Out-Null

# zh-cn, row 3 (d40e0c57...)
$wsZh.Range("F3").Value = $mdUuid2
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $md2Url) | Out-Null

$wsZh.Range("G3").Value = $xlfZh2
$wsZh.Range("G3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $xlfZh2Url) | Out-Null

# de-de, row 2 (0c5cca4e...)
$wsDe.Range("F2").Value = $mdUuid1
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $md1Url) | Out-Null

$wsDe.Range("G2").Value = $xlfDe1
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $xlfDe1Url) | Out-Null

# de-de, row 3 (d40e0c57...)
$wsDe.Range("F3").Value = $mdUuid2
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $md2Url) | Out-Null

$wsDe.Range("G3").Value = $xlfDe2
$wsDe.Range("G3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $xlfDe2Url) | Out-Null
